# credit_score addition + results
# Adds a new column S ("default + 0 heuristic") computed as IF(J>N,1,0),
# extends the S31 average, and adds a new summary row 32 with
# AVERAGE(..2:..24) for columns Q, R and S. Also updates the sheet view
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell S1 -------------------------------------------------
$ws.Range("S1").Value = "default + 0 heuristic"
# Match the look of the other header cells as closely as this host allows
# (copies the "bordered" look already used across row 1 / row 31).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New column S data, rows 2-30 ---------------------------------------
$ws.Range("S2").Formula = "=IF(J2>N2, 1, 0)"
$ws.Range("S3:S30").Formula = "=IF(J3>N3, 1, 0)"

# --- Row 31: average of the new column -----------------------------------
$ws.Range("S31").Formula = "=AVERAGE(S2:S30)"
# Reuse R31's border/fill so the new cell matches the rest of the totals row.
$ws.Range("R31").Copy() | Out-Null
$ws.Range("S31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New row 32: averages restricted to rows 2-24 -------------------------
$ws.Range("Q32").Formula = "=AVERAGE(Q2:Q24)"
$ws.Range("R32").Formula = "=AVERAGE(R2:R24)"
$ws.Range("S32").Formula = "=AVERAGE(S2:S24)"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("Q32:S32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column width for the new column S ------------------------------------
$ws.Columns("S").ColumnWidth = 16.5

# --- Sheet view: scroll / selection ---------------------------------------
$ws.Activate()
$win = $excel.Windows.Item(1)
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("S2").Select()
